# Updates the cryptocurrency price/volume table on Sheet1 (columns D=Price,
# E=Volume(1h)) with refreshed values from the latest data pull.
# Cell contents are plain text (e.g. "26.484.93", "  -0.15%  ") so we force
# a text number format before assigning the value and then restore the
# default "Normal" style, which keeps Excel from re-interpreting strings
# like "242.90" or "0.9975" as numbers (which would drop trailing zeros /
# add floating point noise) while leaving the cell's style index untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.484.93"
Set-TextValue $ws.Range("E2") "  -0.15%  "
Set-TextValue $ws.Range("D3") "1.726.36"
Set-TextValue $ws.Range("E3") "  -0.52%  "
Set-TextValue $ws.Range("E4") "  -0.23%  "
Set-TextValue $ws.Range("D5") "242.90"
Set-TextValue $ws.Range("E5") "  -1.83%  "
Set-TextValue $ws.Range("D6") "0.9975"
Set-TextValue $ws.Range("E6") "  -0.25%  "
Set-TextValue $ws.Range("E7") "  +0.01%  "
Set-TextValue $ws.Range("D8") "0.2605"
Set-TextValue $ws.Range("E8") "  -2.34%  "
Set-TextValue $ws.Range("D9") "0.06199"
Set-TextValue $ws.Range("E9") "  -0.38%  "
Set-TextValue $ws.Range("D10") "1.724.63"
Set-TextValue $ws.Range("E10") "  -0.60%  "
Set-TextValue $ws.Range("D11") "0.06987"
Set-TextValue $ws.Range("E12") "  -0.40%  "
Set-TextValue $ws.Range("D13") "4.529"
Set-TextValue $ws.Range("E13") "  -1.43%  "
Set-TextValue $ws.Range("D14") "0.6013"
Set-TextValue $ws.Range("E14") "  -1.34%  "
Set-TextValue $ws.Range("E16") "  -0.23%  "
Set-TextValue $ws.Range("D17") "26.464.49"
Set-TextValue $ws.Range("D18") "0.9974"
Set-TextValue $ws.Range("E18") "  -0.20%  "
Set-TextValue $ws.Range("D19") "0.000007164"
Set-TextValue $ws.Range("E19") "  -2.03%  "
Set-TextValue $ws.Range("D20") "11.32"
Set-TextValue $ws.Range("E20") "  -1.72%  "
Set-TextValue $ws.Range("D21") "1.948.63"
Set-TextValue $ws.Range("E21") "  -0.35%  "
Set-TextValue $ws.Range("D22") "4.458"
Set-TextValue $ws.Range("E22") "  -2.14%  "
Set-TextValue $ws.Range("D23") "8.516"
Set-TextValue $ws.Range("E23") "  -2.80%  "
Set-TextValue $ws.Range("D24") "5.108"
Set-TextValue $ws.Range("D25") "137.48"
Set-TextValue $ws.Range("E25") "  -2.27%  "
Set-TextValue $ws.Range("D26") "15.30"
Set-TextValue $ws.Range("E26") "  -0.94%  "
Set-TextValue $ws.Range("D27") "1.414"
Set-TextValue $ws.Range("E27") "  -0.89%  "
Set-TextValue $ws.Range("D28") "1.752"
Set-TextValue $ws.Range("E28") "  -1.24%  "
Set-TextValue $ws.Range("D29") "106.58"
Set-TextValue $ws.Range("E29") "  -1.29%  "
Set-TextValue $ws.Range("E30") "  -2.33%  "
Set-TextValue $ws.Range("D31") "0.08009"
Set-TextValue $ws.Range("E31") "  -0.69%  "
Set-TextValue $ws.Range("D32") "3.641"
Set-TextValue $ws.Range("E32") "  -1.35%  "
Set-TextValue $ws.Range("D33") "0.04496"
Set-TextValue $ws.Range("E33") "  -1.44%  "
Set-TextValue $ws.Range("E34") "  -0.26%  "
Set-TextValue $ws.Range("E35") "  -0.51%  "
Set-TextValue $ws.Range("E36") "  -0.49%  "
Set-TextValue $ws.Range("D37") "0.6228"
Set-TextValue $ws.Range("E37") "  -2.01%  "
Set-TextValue $ws.Range("E38") "  +4.00%  "
Set-TextValue $ws.Range("D39") "1.994"
Set-TextValue $ws.Range("E39") "  -1.56%  "
Set-TextValue $ws.Range("D40") "2.385"
Set-TextValue $ws.Range("E40") "  -0.50%  "
Set-TextValue $ws.Range("D41") "0.9971"
Set-TextValue $ws.Range("E41") "  -0.72%  "
Set-TextValue $ws.Range("E42") "  -1.56%  "
Set-TextValue $ws.Range("D43") "99.77"
Set-TextValue $ws.Range("E43") "  -1.71%  "
Set-TextValue $ws.Range("D44") "5.405"
Set-TextValue $ws.Range("E44") "  -0.33%  "
Set-TextValue $ws.Range("D45") "0.3848"
Set-TextValue $ws.Range("E45") "  -1.20%  "
Set-TextValue $ws.Range("D46") "6.899"
Set-TextValue $ws.Range("E46") "  -0.50%  "
Set-TextValue $ws.Range("E47") "  -2.16%  "
Set-TextValue $ws.Range("D48") "0.05371"
Set-TextValue $ws.Range("E48") "  -0.44%  "
Set-TextValue $ws.Range("D49") "30.54"
Set-TextValue $ws.Range("E49") "  +0.06%  "
Set-TextValue $ws.Range("E50") "  -1.25%  "
Set-TextValue $ws.Range("D51") "51.45"
Set-TextValue $ws.Range("E51") "  -0.51%  "
